# Edit script: 
#  1) Swap the three tables' style id from the custom "Table_0" style
#     {E1C2360D-4F5D-417C-B762-D2373FB5EBC8} to the built-in style
#     {BEAE1C6D-6CE9-4744-860C-CFD7D6BA5E34}.
#  2) Re-colour the deck's theme (theme1.xml, the slide master's theme)
#     from the "Integral / Red Violet" palette to the standard
#     "Office Theme / Office" palette.

$p = $ppt.ActivePresentation

# ---- 1) Table style ids -------------------------------------------------
$oldStyleId = "{E1C2360D-4F5D-417C-B762-D2373FB5EBC8}"
$newStyleId = "{BEAE1C6D-6CE9-4744-860C-CFD7D6BA5E34}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            $tbl.ApplyStyle($newStyleId)
        }
    }
}

# ---- 2) Theme colours ----------------------------------------------------
# New (target) "Office" colour scheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink
# Each entry is [R, G, B] (0-255), exactly as the RRGGBB hex reads.
$officeColors = @(
    @(0x00, 0x00, 0x00),  # dk1      000000
    @(0xFF, 0xFF, 0xFF),  # lt1      FFFFFF
    @(0x44, 0x54, 0x6A),  # dk2      44546A
    @(0xE7, 0xE6, 0xE6),  # lt2      E7E6E6
    @(0x5B, 0x9B, 0xD5),  # accent1  5B9BD5
    @(0xED, 0x7D, 0x31),  # accent2  ED7D31
    @(0xA5, 0xA5, 0xA5),  # accent3  A5A5A5
    @(0xFF, 0xC0, 0x00),  # accent4  FFC000
    @(0x44, 0x72, 0xC4),  # accent5  4472C4
    @(0x70, 0xAD, 0x47),  # accent6  70AD47
    @(0x05, 0x63, 0xC1),  # hlink    0563C1
    @(0x95, 0x4F, 0x72)   # folHlink 954F72
)

$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeColors[$i - 1]
    $r = $rgb[0]
    $g = $rgb[1]
    $b = $rgb[2]
    # COM colour packing is 0x00BBGGRR (low byte = R).
    $comRgb = $r -bor ($g -shl 8) -bor ($b -shl 16)
    $tcs.Colors($i).RGB = $comRgb
}
